# The deck's theme is swapped: the slide master's theme (the "Integral" /
# "Red Violet" color scheme) becomes the stock "Office Theme" color scheme
# (and vice-versa for the unused Notes Master theme, which isn't reachable
# through the PowerPoint object model).
#
# ThemeColorScheme.Item(1..12) maps, in order, to:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2,
#   5 accent1, 6 accent2, 7 accent3, 8 accent4, 9 accent5, 10 accent6,
#   11 hlink, 12 folHlink
#
# Target values are the stock Office theme colors (decimal = VBA RGB(r,g,b)
# encoding, i.e. r + g*256 + b*65536):
#   dk1=000000(0) lt1=FFFFFF(16777215) dk2=44546A(6968388) lt2=E7E6E6(15132391)
#   accent1=5B9BD5(13998939) accent2=ED7D31(3243501) accent3=A5A5A5(10855845)
#   accent4=FFC000(49407) accent5=4472C4(12874308) accent6=70AD47(4697456)
#   hlink=0563C1(12673797) folHlink=954F72(7491477)

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$theme = $master.Theme
$cs = $theme.ThemeColorScheme

$cs.Item(1).RGB  = 0         # dk1     000000
$cs.Item(2).RGB  = 16777215  # lt1     FFFFFF
$cs.Item(3).RGB  = 6968388   # dk2     44546A
$cs.Item(4).RGB  = 15132391  # lt2     E7E6E6
$cs.Item(5).RGB  = 13998939  # accent1 5B9BD5
$cs.Item(6).RGB  = 3243501   # accent2 ED7D31
$cs.Item(7).RGB  = 10855845  # accent3 A5A5A5
$cs.Item(8).RGB  = 49407     # accent4 FFC000
$cs.Item(9).RGB  = 12874308  # accent5 4472C4
$cs.Item(10).RGB = 4697456   # accent6 70AD47
$cs.Item(11).RGB = 12673797  # hlink   0563C1
$cs.Item(12).RGB = 7491477   # folHlink 954F72
